$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swap: Cardano (row12) <-> TRON (row13) ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.152"
$ws.Range("E12").Value = "  -4.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "0.364"
$ws.Range("E13").Value = "  +1.54%  "

# --- Simple value updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.478.77"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.600.34"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.81"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.38"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.621.71"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.50"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.063.60"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.497.03"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.615.00"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.37"
$ws.Range("E19").Value = "  +9.09%  "
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.56"
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.01"
$ws.Range("E22").Value = "  +8.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.534"
$ws.Range("E24").Value = "  +13.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.32"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.70"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0786"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  +8.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.39"
$ws.Range("E31").Value = "  +3.06%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.33"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.50"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.22"
$ws.Range("E35").Value = "  +3.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.971"
$ws.Range("E36").Value = "  +9.27%  "
$ws.Range("E37").Value = "  +3.73%  "
$ws.Range("E38").Value = "  +7.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.78"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.83"
$ws.Range("E40").Value = "  +3.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.847"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "295.74"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.27"
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.609"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0984"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.74"
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0547"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.83"
$ws.Range("E51").Value = "  +4.88%  "

# --- Row swap: VeChain (row49) <-> RenderToken (row50) ---
$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "4.93"
$ws.Range("E49").Value = "  +8.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0240"
$ws.Range("E50").Value = "  +2.28%  "

Write-Host "Done applying cryptos update."
